$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores values as plain text in this workbook, e.g.
# "30.156.25" or "1.000" (trailing zeros are significant). Excel normally
# auto-converts such numeric-looking text into a real number when a value is
# assigned, which would corrupt these strings. Apply a Text number format
# first so the assignment below keeps them as text.
$textCells = @('D2','D3','D5','D7','D8','D9','D10','D11','D12','D13','D14','D15','D16','D17','D18','D19','D20','D21','D22','D23','D24','D25','D26','D27','D28','D29','D30','D31','D32','D33','D34','D35','D36','D37','D38','D39','D40','D41','D42','D43','D44','D45','D46','D47','D48','D49','D50','D51')
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = '@'
}

$ws.Range('D2').Value = '30.156.25'
$ws.Range('E2').Value = '  -0.59%  '
$ws.Range('D3').Value = '1.863.42'
$ws.Range('E3').Value = '  -0.43%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '234.06'
$ws.Range('E5').Value = '  -0.93%  '
$ws.Range('E6').Value = '  -0.06%  '
$ws.Range('D7').Value = '0.4676'
$ws.Range('E7').Value = '  -0.82%  '
$ws.Range('B8').Value = 'Cardano'
$ws.Range('C8').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D8').Value = '0.2863'
$ws.Range('E8').Value = '  -1.06%  '
$ws.Range('B9').Value = 'Dogecoin'
$ws.Range('C9').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D9').Value = '0.06490'
$ws.Range('E9').Value = '  -2.30%  '
$ws.Range('B10').Value = 'Solana'
$ws.Range('C10').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D10').Value = '21.15'
$ws.Range('E10').Value = '  -2.53%  '
$ws.Range('B11').Value = 'TRON'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D11').Value = '0.07751'
$ws.Range('E11').Value = '  -3.74%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.862.44'
$ws.Range('E12').Value = '  -0.52%  '
$ws.Range('B13').Value = 'Litecoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D13').Value = '93.82'
$ws.Range('E13').Value = '  -3.76%  '
$ws.Range('B14').Value = 'Polygon'
$ws.Range('C14').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D14').Value = '0.6836'
$ws.Range('E14').Value = '  -0.78%  '
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').Value = '5.050'
$ws.Range('E15').Value = '  -1.82%  '
$ws.Range('B16').Value = 'BitcoinCash'
$ws.Range('C16').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D16').Value = '269.41'
$ws.Range('E16').Value = '  -0.96%  '
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value = '30.159.61'
$ws.Range('E17').Value = '  -0.56%  '
$ws.Range('B18').Value = 'Avalanche'
$ws.Range('C18').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D18').Value = '13.33'
$ws.Range('E18').Value = '  -6.04%  '
$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D19').Value = '0.000007636'
$ws.Range('E19').Value = '  -2.05%  '
$ws.Range('B20').Value = 'Dai'
$ws.Range('C20').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D20').Value = '1.001'
$ws.Range('E20').Value = '  -0.06%  '
$ws.Range('B21').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C21').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D21').Value = '2.125.92'
$ws.Range('E21').Value = '  +0.42%  '
$ws.Range('B22').Value = 'BinanceUSD'
$ws.Range('C22').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D22').Value = '1.000'
$ws.Range('E22').Value = '  -0.15%  '
$ws.Range('B23').Value = 'Uniswap'
$ws.Range('C23').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D23').Value = '5.159'
$ws.Range('E23').Value = '  -3.26%  '
$ws.Range('B24').Value = 'Chainlink'
$ws.Range('C24').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D24').Value = '6.110'
$ws.Range('E24').Value = '  -1.82%  '
$ws.Range('B25').Value = 'Cosmos'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D25').Value = '9.339'
$ws.Range('E25').Value = '  -0.23%  '
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').Value = '165.78'
$ws.Range('E26').Value = '  -1.49%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = '18.55'
$ws.Range('E27').Value = '  -2.42%  '
$ws.Range('B28').Value = 'LidoDAOToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D28').Value = '1.894'
$ws.Range('E28').Value = '  -3.42%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').Value = '1.365'
$ws.Range('E29').Value = '  -0.71%  '
$ws.Range('B30').Value = 'Stellar'
$ws.Range('C30').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D30').Value = '0.09908'
$ws.Range('E30').Value = '  -0.80%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').Value = '1.450'
$ws.Range('E31').Value = '  -1.00%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').Value = '4.226'
$ws.Range('E32').Value = '  -3.57%  '
$ws.Range('B33').Value = 'InternetComputer(DFINITY)'
$ws.Range('C33').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D33').Value = '4.013'
$ws.Range('E33').Value = '  -2.04%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').Value = '0.04686'
$ws.Range('E34').Value = '  -0.57%  '
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').Value = '1.118'
$ws.Range('E35').Value = '  -1.80%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').Value = '0.6887'
$ws.Range('E36').Value = '  -2.06%  '
$ws.Range('B37').Value = 'HuobiToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D37').Value = '2.703'
$ws.Range('E37').Value = '  -0.44%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').Value = '0.01837'
$ws.Range('E38').Value = '  -2.67%  '
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').Value = '2.752'
$ws.Range('E39').Value = '  +3.69%  '
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').Value = '6.348'
$ws.Range('E40').Value = '  +0.34%  '
$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D41').Value = '71.42'
$ws.Range('E41').Value = '  -1.93%  '
$ws.Range('B42').Value = 'PaxDollar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D42').Value = '1.000'
$ws.Range('E42').Value = '  -0.08%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').Value = '1.895'
$ws.Range('E43').Value = '  -3.61%  '
$ws.Range('B44').Value = 'TrustWalletToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D44').Value = '0.8342'
$ws.Range('E44').Value = '  -1.21%  '
$ws.Range('B45').Value = 'Quant'
$ws.Range('C45').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D45').Value = '102.14'
$ws.Range('E45').Value = '  -1.06%  '
$ws.Range('B46').Value = 'TheSandbox'
$ws.Range('C46').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D46').Value = '0.4060'
$ws.Range('E46').Value = '  -2.77%  '
$ws.Range('B47').Value = 'Maker'
$ws.Range('C47').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D47').Value = '936.00'
$ws.Range('E47').Value = '  +0.08%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').Value = '9.101'
$ws.Range('E48').Value = '  -2.16%  '
$ws.Range('B49').Value = 'Aptos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D49').Value = '6.959'
$ws.Range('E49').Value = '  -2.23%  '
$ws.Range('B50').Value = 'Elrond'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D50').Value = '34.05'
$ws.Range('E50').Value = '  -1.53%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').Value = '0.05576'
$ws.Range('E51').Value = '  -1.83%  '

# Revert to the default (Normal) style so the cells don't carry an explicit
# number-format style index that wasn't present in the original file; the
# text value itself is retained.
foreach ($ref in $textCells) {
    $ws.Range($ref).Style = 'Normal'
}
